# Fix a couple of tests with grouping
# The grouped rows in the "Test" sheet contained stray values that broke
# the vertical-grouping expectations. Clear those stray cells and rebuild
# the merged-cell map to match the corrected grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-merge everything first so individual cells inside any merged range
# can be addressed/cleared; merges are rebuilt below from scratch.
$ws.Cells.UnMerge()

# --- Row 3 -----------------------------------------------------------
# B3/D3/E3 were leftover values that don't belong to this group; clear them.
$ws.Range("B3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# --- Row 4 -------------------------------------------------------------
# C4/E4 were leftover values; clear them. F4 becomes an (empty) shared
# string cell instead of the stray date value, same as the blank marker
# cells elsewhere in this sheet (e.g. F6/F7), while keeping its date style.
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F4").Value = "One"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

# --- Row 6 ---------------------------------------------------------------
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F6").ClearContents()

# --- Row 7 ---------------------------------------------------------------
$ws.Range("C7").ClearContents()
$ws.Range("F7").ClearContents()

# --- Row 9 -----------------------------------------------------------------
# C9/D9/E9 were leftover values; clear them. F9 becomes an (empty) shared
# string cell, same treatment as F4 above, keeping its date style.
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("F9").Value = "One"

# --- Rebuild the merged-cell map to match the corrected grouping ---------
$ws.Range("C2:C2").Merge()
$ws.Range("F2:F2").Merge()
$ws.Range("B2:B3").Merge()
$ws.Range("D2:D3").Merge()
$ws.Range("E2:E3").Merge()
$ws.Range("B4:B4").Merge()
$ws.Range("C3:C4").Merge()
$ws.Range("D4:D4").Merge()
$ws.Range("E4:E4").Merge()
$ws.Range("F3:F4").Merge()
$ws.Range("B5:B5").Merge()
$ws.Range("E5:E5").Merge()
$ws.Range("F5:F5").Merge()
$ws.Range("B6:B6").Merge()
$ws.Range("D5:D6").Merge()
$ws.Range("E6:E6").Merge()
$ws.Range("B7:B7").Merge()
$ws.Range("C5:C7").Merge()
$ws.Range("D7:D7").Merge()
$ws.Range("E7:E7").Merge()
$ws.Range("F6:F7").Merge()
$ws.Range("B8:B8").Merge()
$ws.Range("C8:C9").Merge()
$ws.Range("D8:D9").Merge()
$ws.Range("E8:E9").Merge()
$ws.Range("F8:F9").Merge()
